$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Refresh Price (D) / Volume(1h) (E) columns with the latest market data ---
$ws.Range("D2").Value = "'95.565.70"
$ws.Range("E2").Value = "'  +2.47%  "
$ws.Range("D2:E2").ClearFormats()
$ws.Range("D3").Value = "'3.595.80"
$ws.Range("E3").Value = "'  +4.66%  "
$ws.Range("D3:E3").ClearFormats()
$ws.Range("E4").Value = "'  +0.02%  "
$ws.Range("E4").ClearFormats()
$ws.Range("D5").Value = "'238.72"
$ws.Range("E5").Value = "'  +2.77%  "
$ws.Range("D5:E5").ClearFormats()
$ws.Range("D6").Value = "'656.54"
$ws.Range("E6").Value = "'  +5.69%  "
$ws.Range("D6:E6").ClearFormats()
$ws.Range("D7").Value = "'1.48"
$ws.Range("E7").Value = "'  +7.00%  "
$ws.Range("D7:E7").ClearFormats()
$ws.Range("D8").Value = "'0.406"
$ws.Range("E8").Value = "'  +2.99%  "
$ws.Range("D8:E8").ClearFormats()
$ws.Range("E9").Value = "'  -0.09%  "
$ws.Range("E9").ClearFormats()
$ws.Range("E10").Value = "'  +4.05%  "
$ws.Range("E10").ClearFormats()
$ws.Range("D11").Value = "'3.590.13"
$ws.Range("E11").Value = "'  +4.50%  "
$ws.Range("D11:E11").ClearFormats()
$ws.Range("D12").Value = "'43.07"
$ws.Range("E12").Value = "'  +0.72%  "
$ws.Range("D12:E12").ClearFormats()
$ws.Range("E13").Value = "'  +0.86%  "
$ws.Range("E13").ClearFormats()
$ws.Range("D14").Value = "'6.33"
$ws.Range("E14").Value = "'  +1.03%  "
$ws.Range("D14:E14").ClearFormats()
$ws.Range("D15").Value = "'4.262.96"
$ws.Range("E15").Value = "'  +4.66%  "
$ws.Range("D15:E15").ClearFormats()
$ws.Range("D16").Value = "'95.467.66"
$ws.Range("E16").Value = "'  +2.56%  "
$ws.Range("D16:E16").ClearFormats()
$ws.Range("E17").Value = "'  +3.44%  "
$ws.Range("E17").ClearFormats()
$ws.Range("D18").Value = "'3.594.55"
$ws.Range("E18").Value = "'  +4.72%  "
$ws.Range("D18:E18").ClearFormats()
$ws.Range("D19").Value = "'7.91"
$ws.Range("E19").Value = "'  -3.24%  "
$ws.Range("D19:E19").ClearFormats()
$ws.Range("E20").Value = "'  +8.56%  "
$ws.Range("E20").ClearFormats()
$ws.Range("D21").Value = "'17.99"
$ws.Range("E21").Value = "'  -0.77%  "
$ws.Range("D21:E21").ClearFormats()
$ws.Range("D22").Value = "'3.61"
$ws.Range("E22").Value = "'  +7.32%  "
$ws.Range("D22:E22").ClearFormats()
$ws.Range("D23").Value = "'0.497"
$ws.Range("E23").Value = "'  +11.51%  "
$ws.Range("D23:E23").ClearFormats()
$ws.Range("D24").Value = "'509.30"
$ws.Range("E24").Value = "'  +1.19%  "
$ws.Range("D24:E24").ClearFormats()
$ws.Range("E25").Value = "'  +5.80%  "
$ws.Range("E25").ClearFormats()
$ws.Range("D26").Value = "'6.64"
$ws.Range("E26").Value = "'  +0.45%  "
$ws.Range("D26:E26").ClearFormats()
$ws.Range("D27").Value = "'96.86"
$ws.Range("E27").Value = "'  +2.03%  "
$ws.Range("D27:E27").ClearFormats()
$ws.Range("D28").Value = "'12.76"
$ws.Range("E28").Value = "'  +6.75%  "
$ws.Range("D28:E28").ClearFormats()
$ws.Range("D29").Value = "'3.793.67"
$ws.Range("E29").Value = "'  +4.93%  "
$ws.Range("D29:E29").ClearFormats()
$ws.Range("D30").Value = "'3.15"
$ws.Range("E30").Value = "'  +13.37%  "
$ws.Range("D30:E30").ClearFormats()
$ws.Range("E31").Value = "'  -1.25%  "
$ws.Range("E31").ClearFormats()
$ws.Range("E32").Value = "'  +0.12%  "
$ws.Range("E32").ClearFormats()
$ws.Range("E33").Value = "'  +1.75%  "
$ws.Range("E33").ClearFormats()
$ws.Range("E34").Value = "'  +1.15%  "
$ws.Range("E34").ClearFormats()
$ws.Range("E35").Value = "'  +2.41%  "
$ws.Range("E35").ClearFormats()
$ws.Range("D36").Value = "'31.75"
$ws.Range("E36").Value = "'  +4.71%  "
$ws.Range("D36:E36").ClearFormats()
$ws.Range("D37").Value = "'0.561"
$ws.Range("E37").Value = "'  +2.22%  "
$ws.Range("D37:E37").ClearFormats()
$ws.Range("D38").Value = "'574.52"
$ws.Range("E38").Value = "'  +3.06%  "
$ws.Range("D38:E38").ClearFormats()
$ws.Range("E39").Value = "'  +9.24%  "
$ws.Range("E39").ClearFormats()
$ws.Range("E40").Value = "'  +5.27%  "
$ws.Range("E40").ClearFormats()
$ws.Range("E41").Value = "'  -0.01%  "
$ws.Range("E41").ClearFormats()
$ws.Range("E42").Value = "'  +1.14%  "
$ws.Range("E42").ClearFormats()
$ws.Range("D43").Value = "'0.923"
$ws.Range("E43").Value = "'  -0.39%  "
$ws.Range("D43:E43").ClearFormats()
$ws.Range("D44").Value = "'5.73"
$ws.Range("E44").Value = "'  +3.97%  "
$ws.Range("D44:E44").ClearFormats()
$ws.Range("E45").Value = "'  +0.57%  "
$ws.Range("E45").ClearFormats()
$ws.Range("D46").Value = "'23.76"
$ws.Range("E46").Value = "'  +0.33%  "
$ws.Range("D46:E46").ClearFormats()
$ws.Range("D48").Value = "'0.0416"
$ws.Range("E48").Value = "'  +1.46%  "
$ws.Range("D48:E48").ClearFormats()
$ws.Range("D50").Value = "'53.82"
$ws.Range("E50").Value = "'  +0.25%  "
$ws.Range("D50:E50").ClearFormats()
$ws.Range("D51").Value = "'3.49"
$ws.Range("E51").Value = "'  -5.47%  "
$ws.Range("D51:E51").ClearFormats()

# --- Rows 47 & 49: coin ranking swapped places (EnergySwap <-> Stacks) ---
$ws.Range("B47").Value = "Stacks"
$ws.Range("C47").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D47").Value = "'2.25"
$ws.Range("E47").Value = "'  +6.12%  "
$ws.Range("D47:E47").ClearFormats()

$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "'33.62"
$ws.Range("E49").Value = "'  +31.37%  "
$ws.Range("D49:E49").ClearFormats()
